$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title (row 1) ---
$ws.Cells.Item(1, 1).Value = "Relatório de Pagamentos pix do dia: 05/07/2023 - (Todos)"

# --- Remove the 5th data row (old row 8); this shifts the totals row
#     (old row 10) up to row 9, and the table/dimension ranges follow. ---
$ws.Rows.Item(8).Delete()

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value = 57
$ws.Cells.Item(4, 2).Value = "Leonardo Falcão Koblitz"
$ws.Cells.Item(4, 3).Value = 150
$ws.Cells.Item(4, 4).Value = "adm"
$ws.Cells.Item(4, 5).Value = "pago"
$ws.Cells.Item(4, 6).Value = "bra"
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = ""
$ws.Cells.Item(4, 9).Value = ""

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value = 56
$ws.Cells.Item(5, 2).Value = "Eduardo Rossini Xavier da Silva"
$ws.Cells.Item(5, 3).Value = 1934.23
$ws.Cells.Item(5, 4).Value = "adm"
$ws.Cells.Item(5, 5).Value = "pago"
$ws.Cells.Item(5, 6).Value = "adm"
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = ""
$ws.Cells.Item(5, 9).Value = ""

# --- Row 6 ---
$ws.Cells.Item(6, 1).Value = 55
$ws.Cells.Item(6, 2).Value = "Keyce Felix dos Santos"
$ws.Cells.Item(6, 3).Value = 934.23
$ws.Cells.Item(6, 4).Value = "adm"
$ws.Cells.Item(6, 5).Value = "pago"
$ws.Cells.Item(6, 6).Value = "adm"
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = ""

# --- Row 7 ---
$ws.Cells.Item(7, 1).Value = 54
$ws.Cells.Item(7, 2).Value = "Eduardo Rossini Xavier da Silva"
$ws.Cells.Item(7, 3).Value = 924.0700000000001
$ws.Cells.Item(7, 4).Value = "adm"
$ws.Cells.Item(7, 5).Value = "pago"
$ws.Cells.Item(7, 6).Value = "adm"
$ws.Cells.Item(7, 7).Value = ""
$ws.Cells.Item(7, 8).Value = ""
$ws.Cells.Item(7, 9).Value = ""

# --- Totals row (now row 9 after the deletion above) ---
$ws.Cells.Item(9, 2).Value = "Total das taxas de serviço PIX"
$ws.Cells.Item(9, 3).Value = 5
